# Update annotations for Ying Tang
#  - B54 was stored as the text "4"; it should become the number 4.
#  - A new annotation row (55) needs to be appended with the same shape
#    as the existing rows (Annotator, politeness_score, polite_expressions,
#    sentence_purpose, issue_type, id, source_file, text).
#    Note politeness_score on the new row is the text "3" (not numeric),
#    same inline-string style as most other cells in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 54: politeness_score becomes a real number -----------------------
$ws.Range("B54").Value = 4

# --- Row 55: brand-new annotation row --------------------------------------
$ws.Range("A55").Value = "Ying Tang"

# B55 must stay textual ("3"), not become the number 3. Format the cell as
# text before typing the value so Excel doesn't auto-convert it, then put
# the cell style back to Normal so no stray formatting is left behind.
$ws.Range("B55").NumberFormat = "@"
$ws.Range("B55").Value = "3"
$ws.Range("B55").Style = "Normal"

$ws.Range("C55").Value = "We avoid"
$ws.Range("D55").Value = "DIS"
$ws.Range("E55").Value = "MET"
$ws.Range("F55").Value = "0ffe4b07-d72b-4753-8576-ca80ee89bdb3"
$ws.Range("G55").Value = "SJzMATlAZ_annotated.xlsx"
$ws.Range("H55").Value = "We avoid using k-means because it requires knowing the number of clusters a priory."
